$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.109.13'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '3.152.28'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.42'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +2.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '151.27'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +3.86%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '3.150.19'
$ws.Range('E8').Value = '  +2.28%  '
$ws.Range('E9').Value = '  +4.46%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.163'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +6.41%  '
$ws.Range('E11').Value = '  -0.87%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.505'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +7.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000257'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +12.47%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.56'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +6.69%  '
$ws.Range('D15').Value = '3.665.62'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '65.097.51'
$ws.Range('E16').Value = '  +1.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.20'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  +6.52%  '
$ws.Range('D18').Value = '3.146.92'
$ws.Range('E19').Value = '  +0.42%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '510.41'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +6.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.92'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +7.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.734'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +8.58%  '
$ws.Range('E23').Value = '  +13.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.86'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +4.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '85.64'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +5.15%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.93'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +4.50%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.79'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +8.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.20'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +5.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '28.04'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +7.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.999'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  -0.12%  '
$ws.Range('E32').Value = '  +4.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.66'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +6.21%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.04'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +8.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.61'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +6.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '55.67'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '474.75'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +3.41%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0425'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +4.46%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.08'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +1.92%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0857'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +3.69%  '
$ws.Range('D41').Value = '3.116.89'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.63'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +4.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.120'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +3.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.292'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +11.51%  '
$ws.Range('E45').Value = '  +12.82%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '29.26'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +5.02%  '
$ws.Range('D47').Value = '0.0₃0585'
$ws.Range('E47').Value = '  +13.21%  '
$ws.Range('E49').Value = '  +3.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.31'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +11.24%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.63'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -2.34%  '
